$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cell B7 from 2 to 451
$ws.Range("B7").Value = 451

# Add new row 8 (week 7, cases 56)
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 56
